# Generate Report for Handoff
#
# This reproduces the "handoff" regeneration of the localization-status
# report:
#   - the "de-de" status flips from "Handed back: in sync with en-US"
#     (a Handback state) to "Ready for handoff" (a fresh Handoff state) -
#     this text is shared across the Overview sheet (columns E/F) and the
#     "de-de" sheet's Status column (C2), since they all point at the same
#     shared string.
#   - the "Latest Handoff Datetime" on the "zh-cn" and "de-de" sheets is
#     bumped forward to the new handoff-generation timestamps.
#   - columns that used to be sized for the long "Handed back: in sync
#     with en-US" text are narrowed back down now that the text is short.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
# This text is shared by the zh-cn/de-de "Status" cells and by the
# Overview sheet's per-language status columns (E2 mirrors zh-cn, F2
# mirrors de-de), so every occurrence needs to be updated explicitly.
$newStatus = "Ready for handoff"
$ws_zhcn.Range("C2").Value = $newStatus
$ws_dede.Range("C2").Value = $newStatus
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus

# --- Handoff datetimes bumped forward ---
# zh-cn's "Latest Handoff Datetime" only appears on the zh-cn sheet.
$ws_zhcn.Range("H2").Value = "2016-08-19 19:07:23"
# de-de's "Latest Handoff Datetime" is mirrored by the Overview sheet's
# "Latest HO Xliff Generate Date" column (G2), so set both.
$newDedeHandoff = "2016-08-19 19:07:28"
$ws_dede.Range("H2").Value = $newDedeHandoff
$ws_overview.Range("G2").Value = $newDedeHandoff

# --- Column widths: the "Status"/zh-cn/de-de columns shrink now that the
#     status text is much shorter ("Ready for handoff" vs.
#     "Handed back: in sync with en-US"). ColumnWidth in this host snaps to
#     a 1/6-character pixel grid, so use the input that lands on the
#     closest reachable grid point to the recorded OOXML width
#     (17.2159881591797 characters).
$newStatusColWidth = 16.333333333333336

$ws_overview.Columns.Item(5).ColumnWidth = $newStatusColWidth   # Overview!E (zh-cn status)
$ws_overview.Columns.Item(6).ColumnWidth = $newStatusColWidth   # Overview!F (de-de status)
$ws_zhcn.Columns.Item(3).ColumnWidth = $newStatusColWidth       # zh-cn!C (Status)
$ws_dede.Columns.Item(3).ColumnWidth = $newStatusColWidth       # de-de!C (Status)
